$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / date range text update (appears twice: A4 and B7, same shared string) ---
$ws.Range("A4").Value = "1-Jul-25 to 10-Feb-26"
$ws.Range("B7").Value = "1-Jul-25 to 10-Feb-26"

# --- Item description text update ---
$ws.Range("A420").Value = "5617 PATRIKA *-* (M) PUCHNA HAI (7067)"

# --- Numeric data updates ---
# Row 75
$ws.Range("B75").Value = 92
$ws.Range("C75").Value = 179
$ws.Range("E75").Value = 214.8

# Row 91
$ws.Range("B91").Value = 21
$ws.Range("C91").Value = -10
$ws.Range("D91").Value = 0.67
$ws.Range("E91").Value = -6.65

# Row 230
$ws.Range("B230").Value = 51
$ws.Range("C230").Value = 59
$ws.Range("E230").Value = 191.75

# Row 261
$ws.Range("B261").Value = 198
$ws.Range("C261").Value = 99
$ws.Range("E261").Value = 423.72

# Row 367
$ws.Range("B367").Value = 17
$ws.Range("C367").Value = 27.5
$ws.Range("E367").Value = 207.51

# Row 371
$ws.Range("B371").Value = 123
$ws.Range("C371").Value = 90.33
$ws.Range("E371").Value = 406.49

# Row 394
$ws.Range("B394").Value = 16
$ws.Range("C394").Value = 12.5
$ws.Range("E394").Value = 83.13

# Row 409
$ws.Range("B409").Value = 78
$ws.Range("C409").Value = 66.5
$ws.Range("E409").Value = 345.8

# Row 476
$ws.Range("B476").Value = 44
$ws.Range("C476").Value = 9.5
$ws.Range("E476").Value = 98.05

# Row 489
$ws.Range("B489").Value = 18
$ws.Range("C489").Value = 3.5
$ws.Range("E489").Value = 91

# Row 553
$ws.Range("C553").Value = 105.7
$ws.Range("E553").Value = 412.23

# Row 588
$ws.Range("B588").Value = 80
$ws.Range("C588").Value = 210
$ws.Range("E588").Value = 182.7

# Row 589
$ws.Range("B589").Value = 50
$ws.Range("C589").Value = 829
$ws.Range("E589").Value = 480.82

# Row 590
$ws.Range("B590").Value = 66
$ws.Range("C590").Value = 784
$ws.Range("E590").Value = 666.4

# Row 592
$ws.Range("B592").Value = 123
$ws.Range("C592").Value = 219
$ws.Range("E592").Value = 186.15

# Row 658
$ws.Range("B658").Value = 36
$ws.Range("C658").Value = 12.5
$ws.Range("E658").Value = 60

# Row 683
$ws.Range("C683").Value = 22.5
$ws.Range("E683").Value = 95.63

# Row 685
$ws.Range("C685").Value = 38.5
$ws.Range("E685").Value = 192.5

# Row 691
$ws.Range("B691").Value = 10
$ws.Range("C691").Value = 11.5

# Row 729 (Total row)
$ws.Range("C729").Value = 40672.62
$ws.Range("E729").Value = 102442.7
